# "reorg & remove ugly output"
#
# The sheet used to show three columns: the raw metric (A), the model
# name (B) and a recomputed/duplicate value (C) that was only there to
# sanity-check A - and a chunk of rows had been hand-formatted to look
# like a syntax-highlighted code diff (Courier New in grey/orange/green).
# This drops the now-redundant column C output and reorganizes A/B back
# onto one consistent, plain look - matching the style already used by
# the header/edge cells instead of inventing anything new - then updates
# the remembered selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-apply a uniform, plain format (copy the look already used by
#        A1) across the rest of columns A and B ---------------------------
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)
$ws.Range("B2:B9").PasteSpecial(-4122)

# --- 2. Rows 6-7 keep a vertically-centered look (like the header row),
#        just without the old syntax-highlight color ----------------------
$ws.Range("B1").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)

# --- 3. Give column C the same plain look everywhere (including the two
#        rows that never had a C cell before) -----------------------------
$ws.Range("A1").Copy()
$ws.Range("C1:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Remove the "ugly" recomputed output in column C - keep the cells
#        (and their freshly-applied format), just drop the values --------
$ws.Range("C1:C9").ClearContents()

# --- 5. Move the remembered selection -------------------------------------
$ws.Range("B10").Select() | Out-Null
